# Apply cryptocurrency price/volume/name update to Sheet1
# (reflects upstream coinranking.com refresh + new PaxosStandard row,
#  which pushes RenderToken off the bottom of the 50-row list)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.361.73'
$ws.Range("E2").Value = '  -4.00%  '
$ws.Range("E3").Value = '  -3.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4258'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3622'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07040'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8318'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.79%  '
$ws.Range("D12").Value = '1.775.86'
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.220'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.387'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06782'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.64%  '
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008610'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("D21").Value = '26.039.33'
$ws.Range("E21").Value = '  -4.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.992'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.75%  '
$ws.Range("E23").Value = '  +2.35%  '
$ws.Range("D24").Value = '1.954.89'
$ws.Range("E24").Value = '  -5.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.893'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.59%  '
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.024'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.667'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08857'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7212'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.114'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.300'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9996'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.714'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.070'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05082'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01882'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("E40").Value = '  -2.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1600'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.183'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.490'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.81%  '
$ws.Range("B44").Value = 'PaxosStandard'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.69%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.975'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.54%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.24%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06183'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.20%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4470'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.96%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.565'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.56%  '
